$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loot")

# New "№" (index/number) column for the loot sheet - header + running index 0..12
$ws.Range("D1").Value = "№"
$values = 0,1,2,3,4,5,6,7,8,9,10,11,12
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Make "loot" the active sheet/tab and select F12 on it (moves tabSelected
# off of "guns" and onto "loot", and updates the workbook's activeTab)
$ws.Activate() | Out-Null
$ws.Range("F12").Select() | Out-Null

Write-Output "done"
